$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'290.48"
$ws.Range("E2").Value = "'-8.75%"
$ws.Range("D3").Value = "'40.39"
$ws.Range("E3").Value = "'-2.19%"
$ws.Range("D4").Value = "'5.052"
$ws.Range("E4").Value = "'-3.41%"
$ws.Range("D5").Value = "'0.07277"
$ws.Range("E5").Value = "'-5.26%"
$ws.Range("D6").Value = "'4.285"
$ws.Range("E6").Value = "'-1.53%"
$ws.Range("D7").Value = "'1.563"
$ws.Range("E7").Value = "'-7.31%"
$ws.Range("D8").Value = "'0.9187"
$ws.Range("E8").Value = "'-1.77%"
$ws.Range("D9").Value = "'0.1153"
$ws.Range("E9").Value = "'-9.60%"
$ws.Range("E10").Value = "'-6.06%"
$ws.Range("D11").Value = "'0.08660"
$ws.Range("E11").Value = "'-5.25%"
$ws.Range("D12").Value = "'0.04178"
$ws.Range("E12").Value = "'0.84%"
$ws.Range("E13").Value = "'0.42%"
$ws.Range("D14").Value = "'0.001272"
$ws.Range("E14").Value = "'0.27%"
$ws.Range("D15").Value = "'0.005872"
$ws.Range("E15").Value = "'-1.06%"
$ws.Range("D16").Value = "'3.400"
$ws.Range("E16").Value = "'1.68%"
$ws.Range("D18").Value = "'0.3277"
$ws.Range("E18").Value = "'-2.19%"
$ws.Range("D19").Value = "'7.892"
$ws.Range("E19").Value = "'-6.07%"
$ws.Range("E20").Value = "'-1.17%"
$ws.Range("D21").Value = "'0.2886"
$ws.Range("E21").Value = "'0.45%"
$ws.Range("D22").Value = "'0.03868"
$ws.Range("E22").Value = "'-4.27%"
$ws.Range("E23").Value = "'-0.38%"
$ws.Range("D24").Value = "'0.003778"
$ws.Range("E24").Value = "'-7.77%"
$ws.Range("E25").Value = "'0.34%"
$ws.Range("D26").Value = "'0.0003728"
$ws.Range("E38").Value = "'-7.41%"
$ws.Range("D39").Value = "'0.04959"
$ws.Range("E39").Value = "'-6.01%"
$ws.Range("D40").Value = "'0.006431"
$ws.Range("E40").Value = "'212.11%"
$ws.Range("D41").Value = "'0.007687"
$ws.Range("E41").Value = "'-1.51%"
$ws.Range("D42").Value = "'0.1271"
$ws.Range("E42").Value = "'-2.44%"
$ws.Range("D43").Value = "'0.007382"
$ws.Range("E43").Value = "'4.47%"
$ws.Range("D44").Value = "'0.007075"
$ws.Range("E44").Value = "'-15.12%"
$ws.Range("D45").Value = "'0.2900"
$ws.Range("E45").Value = "'-16.35%"
$ws.Range("D46").Value = "'0.00006426"
$ws.Range("E46").Value = "'-3.90%"
$ws.Range("E47").Value = "'-0.48%"
$ws.Range("D48").Value = "'0.01626"
$ws.Range("E48").Value = "'-91.70%"
$ws.Range("E49").Value = "'-0.46%"
$ws.Range("E50").Value = "'-0.48%"
$ws.Range("E51").Value = "'-0.48%"
